$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp text
$ws.Range("A1").Value = "Datos actualizados a 27 de Julio de 2020 a las 00:47"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 4365414
$ws.Range("C4").Value = 49705
$ws.Range("D4").Value = 2085165
$ws.Range("E4").Value = 2130477
$ws.Range("G4").Value = 374
$ws.Range("H4").Value = 149772

# Row 5 - Brasil
$ws.Range("D5").Value = 1634274
$ws.Range("E5").Value = 697813

# Row 23 - Argentina
$ws.Range("B23").Value = 162526
$ws.Range("C23").Value = 4192
$ws.Range("D23").Value = 70518
$ws.Range("E23").Value = 89069
$ws.Range("G23").Value = 46
$ws.Range("H23").Value = 2939

# Row 50 - Nigeria
$ws.Range("B50").Value = 40532
$ws.Range("C50").Value = 555
$ws.Range("D50").Value = 17374
$ws.Range("E50").Value = 22300
$ws.Range("G50").Value = 2
$ws.Range("H50").Value = 858

# Row 51 - Barein
$ws.Range("B51").Value = 39131
$ws.Range("C51").Value = 384
$ws.Range("D51").Value = 35689
$ws.Range("E51").Value = 3302

# Row 59 - Japon
$ws.Range("B59").Value = 29382
$ws.Range("C59").Value = 596
$ws.Range("D59").Value = 21762
$ws.Range("E59").Value = 6624
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = 996

# Row 81 - Bulgaria
$ws.Range("B81").Value = 10427
$ws.Range("C81").Value = 115
$ws.Range("D81").Value = 5355
$ws.Range("E81").Value = 4732
$ws.Range("G81").Value = 2
$ws.Range("H81").Value = 340
